$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing rows with the new combined tuple-style text
$ws.Range("A2").Value = "('Demon', ['Token Creature — Demon', 'Flying', '*/*'])"
$ws.Range("A3").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A4").Value = "('Thrull', ['Token Creature — Thrull', '0/1'])"

# Remove the now-redundant rows that held the individual fields (rows 5-12)
$ws.Range("A5:A12").EntireRow.Delete()
